$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update analysis timestamp cell
$ws.Range("A2").Value = "2025-05-23 14:03:56"

# Update metric values for row 2
$ws.Range("B2").Value = 16219
$ws.Range("C2").Value = 11731
$ws.Range("D2").Value = 72.32875023121031
$ws.Range("E2").Value = 2254
$ws.Range("F2").Value = 13.89728096676737
$ws.Range("H2").Value = 19.15037918490659
$ws.Range("I2").Value = 9552
$ws.Range("J2").Value = 58.89388988223688
$ws.Range("K2").Value = 3075560.17
$ws.Range("L2").Value = 3561
$ws.Range("M2").Value = 21.95573093285653
$ws.Range("N2").Value = 1152511.43
$ws.Range("O2").Value = 4674
$ws.Range("P2").Value = 28.81805290091868
$ws.Range("Q2").Value = 490408.2999999999
$ws.Range("R2").Value = 3545
$ws.Range("S2").Value = 21.85708120106048
$ws.Range("T2").Value = 3447
$ws.Range("U2").Value = 21.25285159380973
$ws.Range("V2").Value = 2436583.87
$ws.Range("W2").Value = 2319
$ws.Range("X2").Value = 14.29804550218879
$ws.Range("Y2").Value = 1431
$ws.Range("Z2").Value = 8.822985387508478
$ws.Range("AA2").Value = 148568
$ws.Range("AB2").Value = 803
$ws.Range("AC2").Value = 4.950983414513842
$ws.Range("AD2").Value = 16231
$ws.Range("AE2").Value = 10610
$ws.Range("AF2").Value = 65.36873883309717
$ws.Range("AG2").Value = 5621
$ws.Range("AH2").Value = 34.63126116690283
$ws.Range("AI2").Value = 470
$ws.Range("AJ2").Value = 868
$ws.Range("AK2").Value = 1365
$ws.Range("AL2").Value = 17.38808731039586
$ws.Range("AM2").Value = 32.11246762856086
$ws.Range("AN2").Value = 50.49944506104328
$ws.Range("AO2").Value = 1375751.25
$ws.Range("AP2").Value = 258284.63
$ws.Range("AQ2").Value = 86070.23
$ws.Range("AR2").Value = 79.98060363845811
$ws.Range("AS2").Value = 15.01562191416203
$ws.Range("AT2").Value = 5.003774447379877
$ws.Range("AU2").Value = 46.96877669165796
$ws.Range("AV2").Value = 241.1653780878748
$ws.Range("AW2").Value = 543.3980938416422
